$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Insert a new row 8 on the "Property1" sheet (pushes the old "Upload"
# row 8 -> 9 and the "Desc" header row 9 -> 10), matching row 7's look.
$ws1.Rows(8).Insert()

$ws1.Cells.Item(8, 1).Value = "Force"
$ws1.Cells.Item(8, 2).Value = $false
$ws1.Cells.Item(8, 3).Value = $false

# Carry the formatting from row 7 onto the freshly inserted row 8.
$ws1.Range("A7:C7").Copy()
$ws1.Range("A8:C8").PasteSpecial(-4122)

# Re-activate Property1 and restore its frozen pane below the new last
# data row (10), then land the selection back on A9 (the "Upload" row).
$ws1.Activate()
$excel.ActiveWindow.FreezePanes = $false
$ws1.Range("A11").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws1.Range("A9").Select()
